$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.013.65'
$ws.Range('E2').Value = '  +3.04%  '
$ws.Range('D3').Value = '3.031.24'
$ws.Range('E3').Value = '  +1.85%  '
$ws.Range('E4').Value = '  +0.10%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '593.62'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -0.10%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '153.61'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +8.21%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.027.72'
$ws.Range('E8').Value = '  +1.84%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.515'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +0.19%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '6.97'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  +17.14%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.152'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +4.27%  '
$ws.Range('E12').Value = '  +2.55%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '35.61'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +4.60%  '
$ws.Range('E15').Value = '  +0.15%  '
$ws.Range('D16').Value = '3.534.29'
$ws.Range('E16').Value = '  +2.05%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '7.09'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +3.75%  '
$ws.Range('D18').Value = '62.974.15'
$ws.Range('E18').Value = '  +2.81%  '
$ws.Range('D19').Value = '3.032.07'
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '453.67'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +1.11%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '14.30'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +1.41%  '
$ws.Range('E22').Value = '  +2.98%  '
$ws.Range('E23').Value = '  +3.90%  '
$ws.Range('B24').Value = 'RenderToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '11.51'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  +12.28%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '83.20'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('E26').Value = '  +9.08%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '12.43'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +4.25%  '
$ws.Range('E29').Value = '  +12.61%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '7.51'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +7.25%  '
$ws.Range('E31').Value = '  +1.68%  '
$ws.Range('E32').Value = '  +0.15%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '27.65'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  +1.94%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '0.111'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('D35').Value = '0.0₃0861'
$ws.Range('E35').Value = '  +6.89%  '
$ws.Range('E36').Value = '  +2.87%  '
$ws.Range('E37').Value = '  +3.09%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '3.14'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  +11.59%  '
$ws.Range('E39').Value = '  +8.78%  '
$ws.Range('E40').Value = '  +3.20%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '50.43'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +1.64%  '
$ws.Range('E43').Value = '  +17.03%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '44.69'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +16.46%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '395.30'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +1.82%  '
$ws.Range('E46').Value = '  +4.00%  '
$ws.Range('D47').Value = '2.720.96'
$ws.Range('E47').Value = '  +1.30%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '133.06'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +2.94%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '25.67'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +11.21%  '
$ws.Range('E50').Value = '  +0.00%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '2.29'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  +8.18%  '
